{"js": "// Edit 1: narrow the \"Ingresos Extraordinarios\" comment down to the\n// \"protecci\u00f3n del menor\" remark.\nconst oldText1 =\n  'El t\u00e9rmino \"Ingresos Extraordinarios\" es un concepto amplio y ambiguo, lo mejor es definir ingresos extraordinarios o que el acuerdo haga referencia a lo que establezca la ley para evitar futuros cambios.';\nconst newText1 =\n  'Para protecci\u00f3n del menor: lo mejor es definir que el acuerdo haga referencia a lo que establezca la ley para evitar futuros cambios.';\n\nconst hits1 = context.document.body.search(oldText1, { matchCase: true });\nhits1.load(\"items\");\nawait context.sync();\nif (hits1.items.length > 0) {\n  hits1.items[0].insertText(newText1, \"Replace\");\n  await context.sync();\n}\n\n// Edit 2: small wording/grammar fix (\"cuando\" -> \"cuanto\", \"que\" -> \"qu\u00e9\").\nconst oldText2 =\n  \"Indican que se debe comprobar previamente los ingresos, indicar cada cuando y de que manera se debe comprobar.\";\nconst newText2 =\n  \"Indican que se debe comprobar previamente los ingresos, indicar cada cuanto y de qu\u00e9 manera se debe comprobar.\";\n\nconst hits2 = context.document.body.search(oldText2, { matchCase: true });\nhits2.load(\"items\");\nawait context.sync();\nif (hits2.items.length > 0) {\n  hits2.items[0].insertText(newText2, \"Replace\");\n  await context.sync();\n}\n\n// Edit 3: drop item \"3 - No hay mecanismo de resoluci\u00f3n de conflictos: Se\n// judicializa r\u00e1pido todo.\" along with its preceding blank paragraph, so\n// \"2 - Educaci\u00f3n...\" flows directly into the blank line before \"4 - Dejar\n// claro...\".\nconst hits3 = context.document.body.search(\n  \"No hay mecanismo de resoluci\u00f3n de conflictos\",\n  { matchCase: true }\n);\nhits3.load(\"items\");\nawait context.sync();\n\nif (hits3.items.length > 0) {\n  const targetParagraph = hits3.items[0].paragraphs.getFirst();\n  const blankParagraph = targetParagraph.getPrevious();\n  targetParagraph.delete();\n  blankParagraph.delete();\n  await context.sync();\n}\n", "ps1": "# Document under edit.\n$d = $word.ActiveDocument\n\n# Edit 1: narrow the \"Ingresos Extraordinarios\" comment down to the\n# \"protecci\u00f3n del menor\" remark.\n$oldText1 = 'El t\u00e9rmino \"Ingresos Extraordinarios\" es un concepto amplio y ambiguo, lo mejor es definir ingresos extraordinarios o que el acuerdo haga referencia a lo que establezca la ley para evitar futuros cambios.'\n$newText1 = 'Para protecci\u00f3n del menor: lo mejor es definir que el acuerdo haga referencia a lo que establezca la ley para evitar futuros cambios.'\n\n$find1 = $d.Content.Find\n$find1.Text = $oldText1\n$find1.Replacement.Text = $newText1\n$find1.Execute($oldText1, $false, $false, $false, $false, $false, $true, 1, $false, $newText1, 2) | Out-Null\n\n# Edit 2: small wording/grammar fix (\"cuando\" -> \"cuanto\", \"que\" -> \"qu\u00e9\").\n$oldText2 = 'Indican que se debe comprobar previamente los ingresos, indicar cada cuando y de que manera se debe comprobar.'\n$newText2 = 'Indican que se debe comprobar previamente los ingresos, indicar cada cuanto y de qu\u00e9 manera se debe comprobar.'\n\n$find2 = $d.Content.Find\n$find2.Text = $oldText2\n$find2.Replacement.Text = $newText2\n$find2.Execute($oldText2, $false, $false, $false, $false, $false, $true, 1, $false, $newText2, 2) | Out-Null\n\n# Edit 3: drop item \"3 - No hay mecanismo de resoluci\u00f3n de conflictos: Se\n# judicializa r\u00e1pido todo.\" along with its preceding blank paragraph, so\n# \"2 - Educaci\u00f3n...\" flows directly into the blank line before \"4 - Dejar\n# claro...\".\n$range3 = $d.Content\n$found3 = $range3.Find.Execute(\"No hay mecanismo de resoluci\u00f3n de conflictos\")\nif ($found3) {\n    $targetPara = $range3.Paragraphs(1)\n    $blankPara = $targetPara.Previous()\n    $targetPara.Range.Delete()\n    $blankPara.Range.Delete()\n}\n"}
